# Update "想去人数" (want-to-go count) values in column F for specific rows
# on both the "展览" and "全部类型" worksheets, which hold identical data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 1898
    7  = 1570
    9  = 620
    13 = 91
    18 = 123
    19 = 3660
    24 = 588
    25 = 304
    28 = 1454
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
